$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - update F column (想去人数 / interest count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13501
$ws1.Range("F4").Value = 656
$ws1.Range("F5").Value = 223
$ws1.Range("F6").Value = 462
$ws1.Range("F7").Value = 1357

# Sheet "全部类型" (all types) - update F column (想去人数 / interest count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13501
$ws4.Range("F4").Value = 656
$ws4.Range("F5").Value = 223
$ws4.Range("F8").Value = 462
$ws4.Range("F9").Value = 1357
